# Update summary table values with the newest airtoxics NATA data.

$wb = $excel.ActiveWorkbook

# --- Sheet "Means": Total Cancer Risk (row 9) and Total Respiratory (row 10) ---
$wsMeans = $wb.Worksheets.Item("Means")

$wsMeans.Range("B9").Value = 26
$wsMeans.Range("C9").Value = 21
$wsMeans.Range("D9").Value = 24
$wsMeans.Range("E9").Value = 21
$wsMeans.Range("F9").Value = 21
$wsMeans.Range("G9").Value = 21

$wsMeans.Range("B10").Value = 0.31
$wsMeans.Range("C10").Value = 0.25
$wsMeans.Range("D10").Value = 0.3
$wsMeans.Range("E10").Value = 0.3
$wsMeans.Range("F10").Value = 0.3
$wsMeans.Range("G10").Value = 0.3

# --- Sheet "Standard Deviations": Total Cancer Risk (row 9) and Total Respiratory (row 10) ---
$wsStdDev = $wb.Worksheets.Item("Standard Deviations")

$wsStdDev.Range("B9").Value = 8.3
$wsStdDev.Range("C9").Value = 3.3
$wsStdDev.Range("D9").Value = 5
$wsStdDev.Range("E9").Value = 3.3
$wsStdDev.Range("F9").Value = 2.4
$wsStdDev.Range("G9").Value = 2.8

$wsStdDev.Range("B10").Value = 0.11
$wsStdDev.Range("C10").Value = 0.063
$wsStdDev.Range("D10").Value = 0
$wsStdDev.Range("E10").Value = 0
$wsStdDev.Range("F10").Value = 0.011
$wsStdDev.Range("G10").Value = 0.023
